$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that (in the source data) held the instructors' names
# with no label in column A. Deleting them shifts every following row up by
# two, which reproduces the row layout (and automatic row-height shift)
# seen in the target workbook.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(13).Delete()

# After the shift, re-point the remaining content cells to the values the
# published workbook ends up with.
$ws.Range("B10").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C10").Value = "8767640 - Eduardo Ferro dos Santos"

$ws.Range("B13").Value = "01/01/2021"
$ws.Range("C13").Value = "01/01/2021"

$ws.Range("B15").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C15").Value = "8767640 - Eduardo Ferro dos Santos"

$ws.Range("B18").Value = "8188658 - Maria Auxiliadora Motta Barreto"
$ws.Range("C18").Value = "8188658 - Maria Auxiliadora Motta Barreto"

$ws.Range("B19").Value = "Aulas Expositivas, Aulas Baseadas em Problemas e Projetos, Atividades Individuais e em Grupo, Seminários. Dadas estas características, haverá múltiplas formas de avaliação definidas pelo docente."
$ws.Range("C19").Value = "Aulas Expositivas, Aulas Baseadas em Problemas e Projetos, Atividades Individuais e em Grupo, Seminários. Dadas estas características, haverá múltiplas formas de avaliação definidas pelo docente."

$ws.Range("B20").Value = "Nota Final = Média Ponderada das formas de avaliação definidas pelo docente, abrangendo avaliações individuais e em grupo."
$ws.Range("C20").Value = "Nota Final = Média Ponderada das formas de avaliação definidas pelo docente, abrangendo avaliações individuais e em grupo."

$ws.Range("B21").Value = "NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota da recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota da recuperação."
